$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.288329
$ws.Range("H2").Value = 117.864987
$ws.Range("I2").Value = 0.632237668435316
$ws.Range("J2").Value = 0.632237668435316
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.398034
$ws.Range("N2").Value = 4.194102
$ws.Range("O2").Value = 0.139066772576779
$ws.Range("P2").Value = 0.139066772576779
$ws.Range("Q2").Value = 54.926419745186
$ws.Range("R2").Value = 494.337777706674
$ws.Range("S2").Value = 0.08792325205076708
$ws.Range("T2").Value = 0.08792325205076709

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.288329
$ws.Range("H3").Value = 117.864987
$ws.Range("I3").Value = 0.632237668435316
$ws.Range("J3").Value = 0.632237668435316
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("N3").Value = 23.817751
$ws.Range("O3").Value = 0.7897418235434783
$ws.Range("P3").Value = 0.7897418235434784
$ws.Range("Q3").Value = 311.9198791093596
$ws.Range("R3").Value = 2807.278911984237
$ws.Range("S3").Value = 0.4993045291829835
$ws.Range("T3").Value = 0.4993045291829836

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.288329
$ws.Range("H4").Value = 117.864987
$ws.Range("I4").Value = 0.632237668435316
$ws.Range("J4").Value = 0.632237668435316
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.715685
$ws.Range("N4").Value = 2.147055
$ws.Range("O4").Value = 0.0711914038797426
$ws.Range("P4").Value = 0.0711914038797426
$ws.Range("Q4").Value = 28.118067740365
$ws.Range("R4").Value = 253.062609663285
$ws.Range("S4").Value = 0.04500988720156537
$ws.Range("T4").Value = 0.04500988720156537

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.344283
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3112930657211948
$ws.Range("J5").Value = 0.3112930657211947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.398034
$ws.Range("N5").Value = 4.194102
$ws.Range("O5").Value = 0.139066772576779
$ws.Range("P5").Value = 0.139066772576779
$ws.Range("Q5").Value = 27.043965339622
$ws.Range("R5").Value = 243.395688056598
$ws.Range("S5").Value = 0.0432905219753777
$ws.Range("T5").Value = 0.0432905219753777

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.344283
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3112930657211948
$ws.Range("J6").Value = 0.3112930657211947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("N6").Value = 23.817751
$ws.Range("O6").Value = 0.7897418235434783
$ws.Range("P6").Value = 0.7897418235434784
$ws.Range("Q6").Value = 153.5791052558443
$ws.Range("R6").Value = 1382.211947302599
$ws.Range("S6").Value = 0.2458411533790962
$ws.Range("T6").Value = 0.2458411533790962

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.344283
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3112930657211948
$ws.Range("J7").Value = 0.3112930657211947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.715685
$ws.Range("N7").Value = 2.147055
$ws.Range("O7").Value = 0.0711914038797426
$ws.Range("P7").Value = 0.0711914038797426
$ws.Range("Q7").Value = 13.844413178855
$ws.Range("R7").Value = 124.599718609695
$ws.Range("S7").Value = 0.02216139036672083
$ws.Range("T7").Value = 0.02216139036672083

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.509096666666667
$ws.Range("H8").Value = 10.52729
$ws.Range("I8").Value = 0.05646926584348937
$ws.Range("J8").Value = 0.05646926584348937
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.398034
$ws.Range("N8").Value = 4.194102
$ws.Range("O8").Value = 0.139066772576779
$ws.Range("P8").Value = 0.139066772576779
$ws.Range("Q8").Value = 4.905836449286666
$ws.Range("R8").Value = 44.15252804358001
$ws.Range("S8").Value = 0.007852998550634208
$ws.Range("T8").Value = 0.00785299855063421

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.509096666666667
$ws.Range("H9").Value = 10.52729
$ws.Range("I9").Value = 0.05646926584348937
$ws.Range("J9").Value = 0.05646926584348937
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("N9").Value = 23.817751
$ws.Range("O9").Value = 0.7897418235434783
$ws.Range("P9").Value = 0.7897418235434784
$ws.Range("Q9").Value = 27.85959688053222
$ws.Range("R9").Value = 250.73637192479
$ws.Range("S9").Value = 0.04459614098139875
$ws.Range("T9").Value = 0.04459614098139875

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.509096666666667
$ws.Range("H10").Value = 10.52729
$ws.Range("I10").Value = 0.05646926584348937
$ws.Range("J10").Value = 0.05646926584348937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.715685
$ws.Range("N10").Value = 2.147055
$ws.Range("O10").Value = 0.0711914038797426
$ws.Range("P10").Value = 0.0711914038797426
$ws.Range("Q10").Value = 2.511407847883333
$ws.Range("R10").Value = 22.60267063095
$ws.Range("S10").Value = 0.004020126311456405
$ws.Range("T10").Value = 0.004020126311456405

Write-Output "Applied Adam15-Itga9 natmi updates"
